# Petty cash book update - 10-Jan-2021, midday update.
# Reset the daily ledger: clear out all entered transactions for the week
# (rows 3-45, columns B:D, plus the first debit formula in D3), set the new
# opening balance and transaction date, and leave the running-balance
# formulas in column E in place (they recompute to the new balance).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New opening balance (SALDO AWAL) on row 2.
$ws.Range("E2").Value = 604525

# Row 3 keeps its "Wages Expense" label but moves to the new date and drops
# its old debit formula (no entries left for the day).
$ws.Range("A3").Value = 44207
$ws.Range("D3").Clear()

# Rows 4-45 had dated transactions (descriptions + debit/credit amounts);
# all of that transaction detail is cleared, leaving just the running
# balance formulas already present in column E.
$ws.Range("B4:D45").Clear()

# Restore the view: scrolled back up with D4 selected.
$ws.Range("D4").Select()
